$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.181.50"
$ws.Range("D3").Value = "1.785.95"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'226.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D8").Value = "'32.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.293"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "'0.0948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "2.044.05"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.790.60"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'11.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "34.156.99"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'67.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "'245.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "'11.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'160.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "'16.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "'0.0520"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'3.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("D34").Value = "'1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").Value = "1.444.39"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("E36").Value = "  +10.57%  "
$ws.Range("D37").Value = "'0.656"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").Value = "'0.0190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").Value = "'81.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").Value = "'13.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "'0.918"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "1.943.78"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "'104.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  -6.26%  "
